$d = $word.ActiveDocument

# Locate the paragraph that ends with "Add line bisection task" (the
# anchor right before the block being rewritten) and the final paragraph
# of the document (which holds the trailing _GoBack bookmark and must be
# left untouched) so the replacement range is computed robustly instead
# of via hard-coded character offsets.
$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Add line bisection task") {
        $startPara = $p
    }
}
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count)

#
# Note: $startPara.Range.End sits exactly at the start of the following
# (empty) paragraph. Because that next paragraph is itself an empty
# placeholder, beginning the replacement range there would leave it
# behind as an untouched paragraph instead of folding it into the
# replacement. Back up one character so the range begins at
# $startPara's own paragraph mark instead -- this consumes that blank
# paragraph into the replacement, exactly as the diff requires.
$replaceStart = $startPara.Range.End - 1
$replaceEnd = $endPara.Range.Start

$r = $d.Range($replaceStart, $replaceEnd)

$quoteL = [char]0x201C
$quoteR = [char]0x201D
$apos = [char]0x2019
$dash = [char]0x2013

$newText = "Remove " + $quoteL + "experimenter" + $quoteR + " text.`r" + `
    "Check fixation cross again. It" + $apos + "s broken.`r" + `
    "Change " + $quoteL + "Ready for next set?" + $quoteR + " to -> TIME FOR A SHORT BREAK (and make it be a break of ~20s).`r" + `
    "Change main sequence to let it be 3 in a row. `r" + `
    "Change practice trials to be 4 for the color. `r" + `
    "Create the gesture part " + $dash + " 2 per format per direction. Gesture for all stimuli. `r" + `
    "`r"

$r.Text = $newText

# The final inserted paragraph (the new blank line right before the
# bookmark paragraph) should carry a left indent of 0.25" (360 twips /
# 18 pt), matching the formatting of the untouched bookmark paragraph
# that follows it. Re-fetch by index (rather than reusing the older
# $endPara reference) since the preceding .Text assignment reshaped the
# paragraph collection.
$newBlankIndex = $d.Paragraphs.Count - 1
$newBlankPara = $d.Paragraphs.Item($newBlankIndex)
$newBlankPara.LeftIndent = 18
